# Update the NATMI Lta-Ltbr ligand-receptor table with freshly recomputed
# TPM-based values, and add the 4 extra rows for the "Resolving-Mac" sending
# cluster (previously only "ECs" was present as a sender).
#
# Each inner array is one data row, in column order A..T:
#   Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#   Ligand-expressing cells, Ligand detection rate,
#   Ligand average expression value, Ligand total expression value,
#   Ligand derived specificity (avg), Ligand derived specificity (total),
#   Receptor-expressing cells, Receptor detection rate,
#   Receptor average expression value, Receptor total expression value,
#   Receptor derived specificity (avg), Receptor derived specificity (total),
#   Edge average expression weight, Edge total expression weight,
#   Edge average expression derived specificity, Edge total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Lta", "Ltbr", "ECs", 2, 0.6666666666666666, 0.4993536666666667, 1.498061, 0.9584860631692095, 0.9584860631692095, 3, 1, 15.01320933333333, 45.039628, 0.2980359842117177, 0.2980359842117177, 7.496901129034223, 67.47211016130801, 0.28566333718985, 0.28566333718985),
    @("ECs", "Lta", "Ltbr", "FAPs", 2, 0.6666666666666666, 0.4993536666666667, 1.498061, 0.9584860631692095, 0.9584860631692095, 3, 1, 17.10605733333334, 51.318172, 0.3395823318071414, 0.3395823318071414, 8.541972451610224, 76.87775206449201, 0.3254849323356472, 0.3254849323356472),
    @("ECs", "Lta", "Ltbr", "MuSCs", 2, 0.6666666666666666, 0.4993536666666667, 1.498061, 0.9584860631692095, 0.9584860631692095, 3, 1, 11.21655833333333, 33.649675, 0.2226664484668797, 0.2226664484668797, 5.601029531130556, 50.40926578017501, 0.2134226875908892, 0.2134226875908892),
    @("ECs", "Lta", "Ltbr", "Resolving-Mac", 2, 0.6666666666666666, 0.4993536666666667, 1.498061, 0.9584860631692095, 0.9584860631692095, 3, 1, 7.037989333333333, 21.113968, 0.1397152355142612, 0.1397152355142612, 3.514445779560889, 31.630012016048, 0.1339151060528231, 0.1339151060528231),
    @("Resolving-Mac", "Lta", "Ltbr", "ECs", 1, 0.3333333333333333, 0.021628, 0.064884, 0.04151393683079058, 0.04151393683079059, 3, 1, 15.01320933333333, 45.039628, 0.2980359842117177, 0.2980359842117177, 0.3247056914613333, 2.922351223152, 0.01237264702186775, 0.01237264702186775),
    @("Resolving-Mac", "Lta", "Ltbr", "FAPs", 1, 0.3333333333333333, 0.021628, 0.064884, 0.04151393683079058, 0.04151393683079059, 3, 1, 17.10605733333334, 51.318172, 0.3395823318071414, 0.3395823318071414, 0.3699698080053334, 3.329728272048, 0.01409739947149423, 0.01409739947149424),
    @("Resolving-Mac", "Lta", "Ltbr", "MuSCs", 1, 0.3333333333333333, 0.021628, 0.064884, 0.04151393683079058, 0.04151393683079059, 3, 1, 11.21655833333333, 33.649675, 0.2226664484668797, 0.2226664484668797, 0.2425917236333333, 2.1833255127, 0.009243760875990532, 0.009243760875990534),
    @("Resolving-Mac", "Lta", "Ltbr", "Resolving-Mac", 1, 0.3333333333333333, 0.021628, 0.064884, 0.04151393683079058, 0.04151393683079059, 3, 1, 7.037989333333333, 21.113968, 0.1397152355142612, 0.1397152355142612, 0.1522176333013333, 1.369958699712, 0.005800129461438068, 0.00580012946143807)
)

$r = 2
foreach ($row in $rows) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r = $r + 1
}
